$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2-6) are being fully replaced, and a new row (7) is being added,
# so clear the existing data block first and rewrite rows 2-7 in one pass.
$ws.Range("A2:AQ6").ClearContents()

$data = New-Object "object[,]" 6,43
$data[0,0] = "Malaysia"
$data[0,1] = 0
$data[0,2] = "Investments & Asset Management"
$data[0,6] = -0.08512690355329949
$data[0,7] = -0.08512690355329949
$data[0,8] = 0.3480710659898477
$data[0,9] = 0.3439349913981692
$data[0,10] = 238.79
$data[0,11] = 12.12131979695432
$data[0,12] = 0
$data[0,13] = 0
$data[0,14] = 0
$data[0,15] = 0
$data[0,16] = 0
$data[0,17] = 0
$data[0,18] = 0
$data[0,20] = 25.604
$data[0,21] = 0.2107793501436533
$data[0,22] = -0.01286549707602339
$data[0,23] = 0.04193426666134865
$data[0,24] = -0.05479976373737205
$data[0,25] = 0.1302419061597149
$data[0,26] = 0.05925756218304797
$data[0,27] = 0.04170090084630979
$data[0,28] = 0.01762175681670111
$data[0,29] = 18.258
$data[0,30] = 0
$data[0,31] = 18.258
$data[0,32] = -7.346000000000004
$data[0,33] = 0.1306653498507847
$data[0,34] = 0.04213494814479763
$data[0,35] = -0.06436688951781791
$data[0,36] = -0.01801735513271429
$data[0,37] = 0.325
$data[0,38] = -0.269
$data[0,39] = 2.545023696682464
$data[0,40] = 21.09846153846154
$data[0,41] = -1.023975466964037
$data[0,42] = -25.49070631970261
$data[1,0] = "Malaysia"
$data[1,1] = "Vertu Capital Limited (LSE:VCBC)"
$data[1,2] = "Investments & Asset Management"
$data[1,10] = -0.17
$data[1,12] = -0
$data[1,13] = -0
$data[1,14] = 0
$data[1,15] = -0
$data[1,16] = -0
$data[1,17] = 0
$data[1,18] = 0
$data[1,20] = 0.314
$data[1,21] = 0.5479930191972078
$data[1,22] = -0.3820224719101124
$data[1,23] = 0.04163324415823549
$data[1,24] = -0.4236557160683479
$data[1,25] = -0
$data[1,26] = 3.695652173913045
$data[1,27] = 0.04163324415823549
$data[1,28] = 3.654018929754809
$data[1,29] = 0
$data[1,30] = 0
$data[1,31] = 0
$data[1,32] = -0.314
$data[1,33] = 0
$data[1,34] = 0
$data[1,35] = -1.212355212355213
$data[1,36] = 6.280000000000001
$data[1,37] = 0
$data[1,38] = 0
$data[2,0] = "Malaysia"
$data[2,1] = "Fintec Global Berhad (KLSE:FINTEC)"
$data[2,2] = "Investments & Asset Management"
$data[2,6] = -0.1725925925925926
$data[2,7] = -0.1725925925925926
$data[2,8] = 0.3925925925925926
$data[2,9] = 0.3925925925925926
$data[2,10] = 240.9
$data[2,11] = 17.84444444444444
$data[2,12] = -0
$data[2,13] = -0
$data[2,14] = -0
$data[2,15] = -0
$data[2,16] = -0
$data[2,17] = -0
$data[2,18] = 0
$data[2,20] = 2.91
$data[2,21] = 0.05069686411149826
$data[2,22] = 4.946611909650924
$data[2,23] = 0.04516939700715081
$data[2,24] = 4.901442512643773
$data[2,25] = 0.2303007557276651
$data[2,26] = 0.09041437076715739
$data[2,27] = 0.04336719649585168
$data[2,28] = 0.0470471742713057
$data[2,29] = 8.26
$data[2,30] = 0
$data[2,31] = 8.26
$data[2,32] = 5.35
$data[2,33] = 0.1257995735607676
$data[2,34] = 0.02596994277809218
$data[2,35] = 0.08525896414342629
$data[2,36] = 0.01697604315405362
$data[2,37] = 0.316
$data[2,38] = 0.213
$data[2,39] = 1.523985239852399
$data[2,40] = 16.77215189873418
$data[2,41] = 0.9870848708487084
$data[2,42] = 24.88262910798122
$data[3,0] = "Malaysia"
$data[3,1] = "OSK Ventures International Berhad (KLSE:OSKVI)"
$data[3,2] = "Investments & Asset Management"
$data[3,6] = 0.4201680672268908
$data[3,7] = 0.4201680672268908
$data[3,8] = 0.680672268907563
$data[3,9] = 0.6402306529306284
$data[3,10] = 2.25
$data[3,11] = 0.6302521008403361
$data[3,12] = -0
$data[3,13] = -0
$data[3,14] = -0
$data[3,15] = -0
$data[3,16] = -0
$data[3,17] = -0
$data[3,18] = 0
$data[3,20] = 1.43
$data[3,21] = 0.05958333333333333
$data[3,22] = 0.05653266331658292
$data[3,23] = 0.04164450690495061
$data[3,24] = 0.0148881564116323
$data[3,25] = 0.09255658396204403
$data[3,26] = 0.05925756218304797
$data[3,27] = 0.04163580536634686
$data[3,28] = 0.01762175681670111
$data[3,29] = 0.011
$data[3,30] = 0
$data[3,31] = 0.011
$data[3,32] = -1.419
$data[3,33] = 0.0004581233601266086
$data[3,34] = 0.0002593666737403032
$data[3,35] = -0.06284044107878305
$data[3,36] = -0.03462580220101998
$data[3,37] = 0
$data[3,38] = 0
$data[3,39] = 0.004489795918367346
$data[3,41] = -0.5791836734693877
$data[4,0] = "Malaysia"
$data[4,1] = "ECM Libra Group Berhad (KLSE:ECM)"
$data[4,2] = "Investments & Asset Management"
$data[4,6] = 0
$data[4,7] = 0
$data[4,8] = 0
$data[4,9] = 0
$data[4,10] = -3.97
$data[4,11] = -1.509505703422053
$data[4,12] = -0
$data[4,13] = -0
$data[4,14] = 0
$data[4,15] = -0
$data[4,16] = -0
$data[4,17] = 0
$data[4,18] = 0
$data[4,20] = 4.15
$data[4,21] = 0.1509090909090909
$data[4,22] = -0.08085539714867618
$data[4,23] = 0.05042600536234269
$data[4,24] = -0.1312814025110189
$data[4,25] = 0.04958521870286575
$data[4,26] = 0
$data[4,27] = 0.04526551957515796
$data[4,28] = -0.04526551957515796
$data[4,29] = 9.84
$data[4,30] = 0
$data[4,31] = 9.84
$data[4,32] = 5.69
$data[4,33] = 0.2635243706480985
$data[4,34] = 0.1774891774891775
$data[4,35] = 0.1714371798734559
$data[4,36] = 0.1109378046402807
$data[4,37] = 0
$data[4,38] = 0
$data[5,0] = "Malaysia"
$data[5,1] = "Pimpinan Ehsan Berhad (KLSE:PEB)"
$data[5,2] = "Investments & Asset Management"
$data[5,10] = -0.22
$data[5,12] = -0
$data[5,13] = -0
$data[5,14] = 0
$data[5,15] = -0
$data[5,16] = -0
$data[5,17] = 0
$data[5,18] = 0
$data[5,20] = 16.8
$data[5,21] = 1.4
$data[5,22] = -0.01286549707602339
$data[5,23] = 0.04193426666134865
$data[5,24] = -0.05479976373737205
$data[5,25] = 0
$data[5,26] = -0.6551724137931032
$data[5,27] = 0.04170090084630979
$data[5,28] = -0.696873314639413
$data[5,29] = 0.147
$data[5,30] = 0
$data[5,31] = 0.147
$data[5,32] = -16.653
$data[5,33] = 0.0121017535193875
$data[5,34] = 0.008572928209016156
$data[5,35] = 3.578981302385556
$data[5,36] = -47.99135446685911
$data[5,37] = 0.008999999999999999
$data[5,38] = -0.482
$data[5,39] = -0.2112068965517241
$data[5,40] = -78.11111111111111
$data[5,41] = 23.92672413793104
$data[5,42] = 1.45850622406639

$ws.Range("A2:AQ7").Value = $data

# The cells below hold numeric-looking text (e.g. a company-name placeholder
# such as "5"); force them to remain text instead of being auto-converted to numbers.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "5"
$ws.Range("B2").Style = "Normal"
